$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.106.56'
$ws.Range("E2").Value = '  -7.27%  '
$ws.Range("D3").Value = '1.419.49'
$ws.Range("E3").Value = '  -7.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9923'
$ws.Range("E4").Value = '  -0.61%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9947'
$ws.Range("E5").Value = '  -0.47%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '273.66'
$ws.Range("E6").Value = '  -5.60%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3689'
$ws.Range("E7").Value = '  -6.38%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3137'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '39.86'
$ws.Range("E9").Value = '  -6.07%  '
$ws.Range("E10").Value = '  -5.25%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06527'
$ws.Range("E11").Value = '  -9.19%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9918'
$ws.Range("E12").Value = '  -0.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.495'
$ws.Range("E13").Value = '  -4.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '17.61'
$ws.Range("E14").Value = '  -4.80%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.182'
$ws.Range("E15").Value = '  -6.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001024'
$ws.Range("E16").Value = '  -6.61%  '
$ws.Range("D17").Value = '1.414.12'
$ws.Range("E17").Value = '  -7.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.05701'
$ws.Range("E18").Value = '  -13.60%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '71.87'
$ws.Range("E19").Value = '  -14.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9932'
$ws.Range("E20").Value = '  -0.59%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.600'
$ws.Range("E21").Value = '  -9.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.93'
$ws.Range("E22").Value = '  -4.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.02'
$ws.Range("E23").Value = '  +2.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.253'
$ws.Range("E24").Value = '  -5.16%  '
$ws.Range("D25").Value = '20.107.82'
$ws.Range("E25").Value = '  -7.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.283'
$ws.Range("E26").Value = '  -4.08%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '136.50'
$ws.Range("E27").Value = '  -9.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '16.99'
$ws.Range("E28").Value = '  -7.91%  '
$ws.Range("D29").Value = '1.569.69'
$ws.Range("E29").Value = '  -8.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '109.72'
$ws.Range("E30").Value = '  -6.67%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.146'
$ws.Range("E31").Value = '  -14.79%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.348'
$ws.Range("E32").Value = '  -12.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.8366'
$ws.Range("E33").Value = '  -13.15%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.07684'
$ws.Range("E34").Value = '  -5.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '8.460'
$ws.Range("E35").Value = '  -1.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.480'
$ws.Range("E36").Value = '  +0.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05839'
$ws.Range("E37").Value = '  -3.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.854'
$ws.Range("E38").Value = '  -6.98%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.9937'
$ws.Range("E39").Value = '  -0.51%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02084'
$ws.Range("E40").Value = '  -6.81%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '10.56'
$ws.Range("E41").Value = '  -5.70%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1910'
$ws.Range("E42").Value = '  -6.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.094'
$ws.Range("E43").Value = '  -7.82%  '
$ws.Range("E44").Value = '  -9.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.24'
$ws.Range("E45").Value = '  -6.61%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.517'
$ws.Range("E46").Value = '  -5.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5182'
$ws.Range("E47").Value = '  -7.28%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '112.14'
$ws.Range("E48").Value = '  -4.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.780'
$ws.Range("E49").Value = '  -5.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.038'
$ws.Range("E50").Value = '  -11.28%  '

$ws.Range("B51").Value = 'PaxDollar'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.9919'
$ws.Range("E51").Value = '  -0.63%  '
